$wb = $excel.ActiveWorkbook

# --- optimization_parameters sheet (the bulk of the edit) ---
$ws7 = $wb.Worksheets.Item("optimization_parameters")

# Row 1: keep only A1:B1, drop the extra (stray) "value" cells in C1:F1
$ws7.Range("C1:F1").ClearContents()

# Insert a new row 9 for the "L_curve" parameter (rows 9-16 shift down to 10-17)
$ws7.Rows.Item(9).Insert()

# A8's label changed from "Model" to "production_function"
$ws7.Range("A8").Value = "production_function"

# Populate the freshly-inserted row 9
$ws7.Range("A9").Value = "L_curve"
$ws7.Range("B9").Value = 0
$ws7.Range("B9").NumberFormat = $ws7.Range("B2").NumberFormat

# The old "Deletion" row (now at row 17, after the insert above) is removed entirely
$ws7.Rows.Item(17).Delete()

# Update the sheet's recorded selection
$ws7.Range("C1:F4").Select()

# --- network_weights sheet selection stays put; just the "active tab" marker moves ---
$ws6 = $wb.Worksheets.Item("network_weights")
$ws6.Range("D8").Select()

# optimization_parameters becomes the active tab
$ws7.Activate()
